# 12/12/14 hours update/descriptions added
#
# This script reproduces the commit: updates a few existing hours entries on
# the "Work Database" sheet, appends two new timesheet rows for 12/12/2014,
# and adds the missing Sub-Task descriptions on the "Sub-Tasks" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Work Database" sheet - correct a couple of existing entries and log the
#    hours worked on 12/12/2014 (two new rows at the bottom of the table).
# ---------------------------------------------------------------------------
$wsDb = $wb.Worksheets.Item("Work Database")
$wsDb.Activate()

# Row 14 was logged under the wrong month.
$wsDb.Range("C14").Value = "November"

# A few existing hour totals needed correcting.
$wsDb.Range("E17").Value = 5
$wsDb.Range("E18").Value = 6.5
$wsDb.Range("E20").Value = 4

# New 12/12/2014 entries.
$wsDb.Range("A47").Value = "Reporting  - Git Hub"
$wsDb.Range("B47").Value = "Update"
$wsDb.Range("C47").Value = "December"
$wsDb.Range("D47").Value = 2014
$wsDb.Range("E47").Value = 3
$wsDb.Range("F47").Value = "April"

$wsDb.Range("A48").Value = "LCD Panels"
$wsDb.Range("B48").Value = "Update"
$wsDb.Range("C48").Value = "December"
$wsDb.Range("D48").Value = 2014
$wsDb.Range("E48").Value = 0.5
$wsDb.Range("F48").Value = "April"

$wsDb.Range("I26").Select()

# ---------------------------------------------------------------------------
# 2. "Sub-Tasks" sheet - fill in the descriptions that were still missing.
# ---------------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Sub-Tasks")
$wsSub.Activate()

$wsSub.Range("B20").Value = "Going through the motions of getting Ensemble to do what I need it to do for various projects"
$wsSub.Range("B23").Value = "Communication regarding the task at hand"
$wsSub.Range("B24").Value = "University required training"
$wsSub.Range("B26").Value = "Answering general questions not related to any projects"
$wsSub.Range("B27").Value = "Making sure the TRB demo site is live, functioning and accessible to those who need it"
$wsSub.Range("B28").Value = "Maintaining the hardware (mostly the computers for LCD panels), making sure they're updating and not displaying errors"

$wsSub.Range("B6").Select()

# Leave "Work Database" as the active sheet/tab.
$wsDb.Activate()
